$wb = $excel.ActiveWorkbook

# --- groups sheet: selection moves from E18 to L34 ---
$wsGroups = $wb.Worksheets.Item("groups")
$wsGroups.Range("L34").Select()

# --- TechColor sheet: rename "ENERGY_INTERCHANGE" row to the new
#     "ENERGY_IMPORT_S2" / description, highlight it, and move the view ---
$wsTechColor = $wb.Worksheets.Item("TechColor")
$wsTechColor.Range("D67").Value = "Import of Electricity out of the sate"
$wsTechColor.Range("A67").Value = "ENERGY_IMPORT_S2"
$wsTechColor.Range("A67").Interior.Color = 65535
$wsTechColor.Range("D62").Select()

# --- SummaryGroupsExistingTech sheet: selection moves from C14 to F31 ---
$wsSummary = $wb.Worksheets.Item("SummaryGroupsExistingTech")
$wsSummary.Range("F31").Select()

# --- KeepTechList sheet: split the "ENERGY_INTERCHANGE" row into two new
#     rows ("ENERGY_IMPORT_S1" and "ENERGY_IMPORT_S2"), pushing the
#     remaining rows down, and highlight the two new rows ---
$wsKeepTech = $wb.Worksheets.Item("KeepTechList")
$wsKeepTech.Rows.Item(59).Insert()
$wsKeepTech.Range("A58").Value = "ENERGY_IMPORT_S1"
$wsKeepTech.Range("A59").Value = "ENERGY_IMPORT_S2"
$wsKeepTech.Range("A58").Interior.Color = 65535
$wsKeepTech.Range("A59").Interior.Color = 65535
$wsKeepTech.Range("E59").Select()

# --- MaxCapacityGroup sheet becomes the active tab (selection stays B2) ---
$wsMaxCap = $wb.Worksheets.Item("MaxCapacityGroup")
$wsMaxCap.Activate()
